# Fixing and updating to new digisales mobile
# The "URL_DIGISALESMOBILE" row's URL (cell E2 on the "Global" sheet) is
# being updated to point at the new Digisales Mobile test server.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")
$ws.Activate()

# Update the Digisales Mobile URL from the old host/port to the new one
$ws.Range("E2").Value = "http://192.168.231.13:99/"

# Match the author's final cursor position / scroll state: E2 selected
# and the sheet scrolled back to show column A (no more topLeftCell override)
[void]$ws.Range("E2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
